$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.628.63"
$ws.Range("E2").Value = "  +3.98%  "

$ws.Range("D3").Value = "3.070.65"
$ws.Range("E3").Value = "  +2.64%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.26%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.064.08"
$ws.Range("E8").Value = "  +2.47%  "

$ws.Range("E9").Value = "  +1.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.32%  "

$ws.Range("E11").Value = "  +2.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000229"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.04%  "

$ws.Range("D15").Value = "3.564.68"
$ws.Range("E15").Value = "  +2.76%  "

$ws.Range("D16").Value = "63.620.75"
$ws.Range("E16").Value = "  +3.95%  "

$ws.Range("D17").Value = "3.067.89"
$ws.Range("E17").Value = "  +2.78%  "

$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "485.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.678"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.71%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.03%  "

$ws.Range("E29").Value = "  +6.85%  "

$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.97%  "

$ws.Range("E32").Value = "  +0.32%  "

$ws.Range("E33").Value = "  +7.93%  "

$ws.Range("E34").Value = "  +5.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "469.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0827"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0399"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.91%  "

$ws.Range("D40").Value = "3.062.65"
$ws.Range("E40").Value = "  -3.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.120"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "28.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.257"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.49%  "

$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.20%  "

$ws.Range("E48").Value = "  +2.05%  "

$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").Value = "0.0₃0516"
$ws.Range("E49").Value = "  +4.57%  "

$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.72%  "

$ws.Range("E51").Value = "  +4.34%  "
